$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("2016")
$ws.PageSetup.Orientation = 1
$ws.PageSetup.PrintQuality = 200
Write-Host "done"
